# Update DAMSLTag (column I) and DialogAct (column J) values for specific rows
# following a re-run of SGNN dialog-act annotation after transcript clean up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 5; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 11; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 18; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 32; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 36; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 41; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 46; I = 'ba'; J = 'Appreciation' },
    @{ Row = 56; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 57; I = 'ba'; J = 'Appreciation' },
    @{ Row = 62; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 86; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 87; I = 'ba'; J = 'Appreciation' },
    @{ Row = 94; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 101; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 108; I = 'ba'; J = 'Appreciation' },
    @{ Row = 109; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 110; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 111; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 113; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 116; I = '%'; J = 'Uninterpretable' },
    @{ Row = 144; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 149; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 163; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 164; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 168; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 171; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 178; I = 'ba'; J = 'Appreciation' },
    @{ Row = 179; I = 'ba'; J = 'Appreciation' },
    @{ Row = 183; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 185; I = 'ba'; J = 'Appreciation' },
    @{ Row = 191; I = 'ba'; J = 'Appreciation' },
    @{ Row = 195; I = 'ba'; J = 'Appreciation' },
    @{ Row = 202; I = 'ba'; J = 'Appreciation' },
    @{ Row = 226; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 230; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 231; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 234; I = 'ba'; J = 'Appreciation' },
    @{ Row = 254; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 296; I = '%'; J = 'Uninterpretable' },
    @{ Row = 299; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 300; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 319; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 324; I = 'ba'; J = 'Appreciation' },
    @{ Row = 351; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 353; I = 'ba'; J = 'Appreciation' },
    @{ Row = 359; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 384; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 408; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 420; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 430; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 438; I = 'ba'; J = 'Appreciation' },
    @{ Row = 439; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 441; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 445; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 449; I = 'ba'; J = 'Appreciation' },
    @{ Row = 456; I = 'ba'; J = 'Appreciation' },
    @{ Row = 458; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 459; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 465; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 472; I = '%'; J = 'Uninterpretable' },
    @{ Row = 479; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 484; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 498; I = 'ba'; J = 'Appreciation' },
    @{ Row = 499; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 500; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 527; I = 'ba'; J = 'Appreciation' },
    @{ Row = 534; I = 'sv'; J = 'Statement-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
